# Updates the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# ranking sheet with refreshed values, matching a scheduled scrape run.
#
# Several "Price" values read like plain decimals (e.g. "1.014") which
# Excel's Range.Value setter would otherwise auto-convert to a Double,
# losing the original text formatting. A leading apostrophe forces those
# assignments to stay text (Excel's normal "quote prefix" behaviour),
# while values that already fail numeric parsing (e.g. "27.506.55", which
# has two dots) are assigned as plain text without needing the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.506.55"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.868.42"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'311.96"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4778"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.3744"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "'0.9347"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "'20.67"
$ws.Range("E11").Value = "  +5.01%  "
$ws.Range("D12").Value = "'0.07828"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.882.54"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "'5.440"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'6.555"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'90.15"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'0.000008889"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "27.538.18"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "'14.60"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'5.116"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'1.950"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "'154.49"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'2.024"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "'115.66"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "'0.08901"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'3.336"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "'1.220"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("D33").Value = "'0.7595"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'2.732"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'0.02036"
$ws.Range("E36").Value = "  +4.00%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'2.997"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "'0.05268"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.5308"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "'8.479"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").Value = "'10.59"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "'1.013"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "'1.654"
$ws.Range("E47").Value = "  +3.15%  "
$ws.Range("D48").Value = "'102.81"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "'67.35"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "'0.06080"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "'0.9201"
$ws.Range("E51").Value = "  +3.80%  "
